# Flip the "N" sensitivity flags to "Y" across the four variable sheets and
# move each sheet's active-cell selection to reflect where editing stopped.

$wb = $excel.ActiveWorkbook

# --- Globals ---
$ws = $wb.Worksheets.Item("Globals")
$ws.Range("B5").Value = "Y"
$ws.Range("B5").Select()

# --- PowerPlants ---
$ws = $wb.Worksheets.Item("PowerPlants")
$ws.Range("B2:B8").Value = "Y"
$ws.Range("B11").Select()

# --- Fuels ---
$ws = $wb.Worksheets.Item("Fuels")
$ws.Range("B2:B6").Value = "Y"
$ws.Range("B7").Select()

# --- Connections ---
$ws = $wb.Worksheets.Item("Connections")
$ws.Range("B2:B5").Value = "Y"

# Leave the Globals sheet active/selected, matching tabSelected="1" in the source.
$wb.Worksheets.Item("Globals").Activate()
